# Generate Report for Handback
# Simulate a failed handback report generation: clear the "Latest Target File"
# and "Latest Handback File" columns, stamp the "Latest Handback DateTime"
# with the default/error datetime, and record the error message, for both
# the zh-cn and de-de sheets. Then let the affected columns auto-fit.

$wb = $excel.ActiveWorkbook

$sheetNames = @("zh-cn", "de-de")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Remove the "Latest Target File" hyperlink (column I) - handback failed
    # before a target file could be produced.
    $ws.Hyperlinks.Item(1).Delete()

    # Clear Latest Target File (I2) and Latest Handback File (J2)
    $ws.Range("I2").Value = ""
    $ws.Range("J2").Value = ""

    # Stamp Latest Handback DateTime (K2) with the default/error datetime
    $ws.Range("K2").Value = "0001-01-01 00:00:00"

    # Record the error detail (P2)
    $ws.Range("P2").Value = "The given key was not present in the dictionary."

    # Auto-fit the columns whose content changed
    $ws.Range("I:J").EntireColumn.AutoFit() | Out-Null
    $ws.Range("P:P").EntireColumn.AutoFit() | Out-Null
}
